{"js": "// Auto-generated replacements: sequential find & replace, order matters\nconst pairs = [\n  [\"2023-07-13 Thursday\", \"2023-07-14 Friday\"],\n  [\"44\u00d760=2640\", \"52\u00d729=1508\"],\n  [\"70\u00d745=3150\", \"18\u00d775=1350\"],\n  [\"32\u00d785=2720\", \"89\u00d724=2136\"],\n  [\"87\u00d733=2871\", \"100\u00d796=9600\"],\n  [\"64\u00d738=2432\", \"60\u00d792=5520\"],\n  [\"60\u00d795=5700\", \"14\u00d761=854\"],\n  [\"15\u00d779=1185\", \"77\u00d765=5005\"],\n  [\"17\u00d721=357\", \"38\u00d717=646\"],\n  [\"34\u00d721=714\", \"15\u00d740=600\"],\n  [\"35\u00d747=1645\", \"73\u00d786=6278\"],\n  [\"13\u00d740=520\", \"81\u00d768=5508\"],\n  [\"20\u00d732=640\", \"71\u00d772=5112\"],\n  [\"53\u00d780=4240\", \"94\u00d766=6204\"],\n  [\"37\u00d794=3478\", \"37\u00d763=2331\"],\n  [\"87\u00d773=6351\", \"50\u00d770=3500\"],\n  [\"22\u00d776=1672\", \"65\u00d760=3900\"],\n  [\"56\u00d772=4032\", \"38\u00d798=3724\"],\n  [\"76\u00d774=5624\", \"81\u00d721=1701\"],\n  [\"69\u00d728=1932\", \"56\u00d759=3304\"],\n  [\"70\u00d759=4130\", \"95\u00d785=8075\"],\n  [\"43\u00d767=2881\", \"31\u00d724=744\"],\n  [\"73\u00d753=3869\", \"44\u00d755=2420\"],\n  [\"13\u00d744=572\", \"57\u00d748=2736\"],\n  [\"23\u00d744=1012\", \"33\u00d789=2937\"],\n  [\"10\u00d713=130\", \"98\u00d768=6664\"],\n  [\"98\u00d757=5586\", \"18\u00d763=1134\"],\n  [\"83\u00d756=4648\", \"11\u00d710=110\"],\n  [\"17\u00d732=544\", \"52\u00d755=2860\"],\n  [\"79\u00d776=6004\", \"96\u00d769=6624\"],\n  [\"11\u00d786=946\", \"68\u00d738=2584\"],\n  [\"43\u00d780=3440\", \"21\u00d738=798\"],\n  [\"64\u00d728=1792\", \"43\u00d781=3483\"],\n  [\"32\u00d758=1856\", \"67\u00d761=4087\"],\n  [\"92\u00d774=6808\", \"26\u00d727=702\"],\n  [\"12\u00d771=852\", \"60\u00d761=3660\"],\n  [\"23\u00d781=1863\", \"27\u00d753=1431\"],\n  [\"50\u00d788=4400\", \"13\u00d790=1170\"],\n  [\"74\u00d769=5106\", \"21\u00d776=1596\"],\n  [\"76\u00d713=988\", \"60\u00d726=1560\"],\n  [\"27\u00d720=540\", \"92\u00d733=3036\"],\n  [\"31\u00d712=372\", \"44\u00d752=2288\"],\n  [\"89\u00d765=5785\", \"97\u00d757=5529\"],\n  [\"33\u00d740=1320\", \"97\u00d719=1843\"],\n  [\"50\u00d721=1050\", \"99\u00d723=2277\"],\n  [\"71\u00d754=3834\", \"54\u00d713=702\"],\n  [\"78\u00d716=1248\", \"82\u00d775=6150\"],\n  [\"44\u00d770=3080\", \"24\u00d794=2256\"],\n  [\"83\u00d776=6308\", \"11\u00d758=638\"],\n  [\"63\u00d774=4662\", \"25\u00d788=2200\"],\n  [\"97\u00d723=2231\", \"10\u00d796=960\"],\n  [\"47\u00d762=2914\", \"91\u00d789=8099\"],\n  [\"99\u00d780=7920\", \"62\u00d794=5828\"],\n  [\"57\u00d790=5130\", \"30\u00d742=1260\"],\n  [\"10\u00d746=460\", \"11\u00d756=616\"],\n  [\"90\u00d723=2070\", \"71\u00d765=4615\"],\n  [\"16\u00d727=432\", \"17\u00d728=476\"],\n  [\"52\u00d739=2028\", \"59\u00d766=3894\"],\n  [\"87\u00d752=4524\", \"29\u00d722=638\"],\n  [\"15\u00d771=1065\", \"56\u00d747=2632\"],\n  [\"37\u00d712=444\", \"96\u00d775=7200\"],\n  [\"90\u00d740=3600\", \"57\u00d743=2451\"],\n  [\"42\u00d743=1806\", \"16\u00d773=1168\"],\n  [\"30\u00d760=1800\", \"19\u00d757=1083\"],\n  [\"84\u00d775=6300\", \"88\u00d735=3080\"],\n  [\"22\u00d761=1342\", \"62\u00d744=2728\"],\n  [\"16\u00d717=272\", \"21\u00d775=1575\"],\n  [\"83\u00d777=6391\", \"51\u00d765=3315\"],\n  [\"78\u00d764=4992\", \"13\u00d778=1014\"],\n  [\"24\u00d748=1152\", \"22\u00d765=1430\"],\n  [\"100\u00d735=3500\", \"46\u00d739=1794\"],\n  [\"14\u00d740=560\", \"15\u00d761=915\"],\n  [\"84\u00d798=8232\", \"20\u00d781=1620\"],\n  [\"45\u00d729=1305\", \"65\u00d734=2210\"],\n  [\"85\u00d715=1275\", \"13\u00d719=247\"],\n  [\"72\u00d780=5760\", \"45\u00d736=1620\"],\n  [\"53\u00d729=1537\", \"57\u00d773=4161\"],\n  [\"92\u00d742=3864\", \"35\u00d775=2625\"],\n  [\"19\u00d775=1425\", \"93\u00d729=2697\"],\n  [\"36\u00d749=1764\", \"11\u00d751=561\"],\n  [\"54\u00d733=1782\", \"22\u00d749=1078\"],\n  [\"12\u00d773=876\", \"28\u00d778=2184\"],\n  [\"14\u00d760=840\", \"18\u00d729=522\"],\n  [\"99\u00d733=3267\", \"78\u00d786=6708\"],\n  [\"58\u00d776=4408\", \"20\u00d760=1200\"],\n  [\"92\u00d773=6716\", \"45\u00d768=3060\"],\n  [\"98\u00d745=4410\", \"78\u00d740=3120\"],\n  [\"66\u00d730=1980\", \"58\u00d759=3422\"],\n  [\"22\u00d726=572\", \"52\u00d711=572\"],\n  [\"70\u00d785=5950\", \"16\u00d750=800\"],\n  [\"97\u00d761=5917\", \"55\u00d713=715\"],\n  [\"94\u00d756=5264\", \"27\u00d727=729\"],\n  [\"62\u00d726=1612\", \"87\u00d766=5742\"],\n  [\"56\u00d790=5040\", \"45\u00d773=3285\"],\n  [\"77\u00d735=2695\", \"85\u00d720=1700\"],\n  [\"54\u00d765=3510\", \"14\u00d755=770\"],\n  [\"15\u00d787=1305\", \"97\u00d759=5723\"],\n  [\"86\u00d767=5762\", \"28\u00d718=504\"],\n  [\"40\u00d789=3560\", \"64\u00d786=5504\"],\n  [\"96\u00d748=4608\", \"17\u00d772=1224\"],\n  [\"88\u00d717=1496\", \"36\u00d749=1764\"],\n];\n\n// Apply each replacement in order. Every \"old\" value below is unique in the\n// document at the moment it is processed (verified against the source diff),\n// so a plain search + replace-first-match is unambiguous and order-safe even\n// though a handful of \"new\" values coincide with \"old\" values used earlier.\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Auto-generated replacements: sequential find & replace, order matters\n$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@('2023-07-13 Thursday', '2023-07-14 Friday')\n    ,@('44\u00d760=2640', '52\u00d729=1508')\n    ,@('70\u00d745=3150', '18\u00d775=1350')\n    ,@('32\u00d785=2720', '89\u00d724=2136')\n    ,@('87\u00d733=2871', '100\u00d796=9600')\n    ,@('64\u00d738=2432', '60\u00d792=5520')\n    ,@('60\u00d795=5700', '14\u00d761=854')\n    ,@('15\u00d779=1185', '77\u00d765=5005')\n    ,@('17\u00d721=357', '38\u00d717=646')\n    ,@('34\u00d721=714', '15\u00d740=600')\n    ,@('35\u00d747=1645', '73\u00d786=6278')\n    ,@('13\u00d740=520', '81\u00d768=5508')\n    ,@('20\u00d732=640', '71\u00d772=5112')\n    ,@('53\u00d780=4240', '94\u00d766=6204')\n    ,@('37\u00d794=3478', '37\u00d763=2331')\n    ,@('87\u00d773=6351', '50\u00d770=3500')\n    ,@('22\u00d776=1672', '65\u00d760=3900')\n    ,@('56\u00d772=4032', '38\u00d798=3724')\n    ,@('76\u00d774=5624', '81\u00d721=1701')\n    ,@('69\u00d728=1932', '56\u00d759=3304')\n    ,@('70\u00d759=4130', '95\u00d785=8075')\n    ,@('43\u00d767=2881', '31\u00d724=744')\n    ,@('73\u00d753=3869', '44\u00d755=2420')\n    ,@('13\u00d744=572', '57\u00d748=2736')\n    ,@('23\u00d744=1012', '33\u00d789=2937')\n    ,@('10\u00d713=130', '98\u00d768=6664')\n    ,@('98\u00d757=5586', '18\u00d763=1134')\n    ,@('83\u00d756=4648', '11\u00d710=110')\n    ,@('17\u00d732=544', '52\u00d755=2860')\n    ,@('79\u00d776=6004', '96\u00d769=6624')\n    ,@('11\u00d786=946', '68\u00d738=2584')\n    ,@('43\u00d780=3440', '21\u00d738=798')\n    ,@('64\u00d728=1792', '43\u00d781=3483')\n    ,@('32\u00d758=1856', '67\u00d761=4087')\n    ,@('92\u00d774=6808', '26\u00d727=702')\n    ,@('12\u00d771=852', '60\u00d761=3660')\n    ,@('23\u00d781=1863', '27\u00d753=1431')\n    ,@('50\u00d788=4400', '13\u00d790=1170')\n    ,@('74\u00d769=5106', '21\u00d776=1596')\n    ,@('76\u00d713=988', '60\u00d726=1560')\n    ,@('27\u00d720=540', '92\u00d733=3036')\n    ,@('31\u00d712=372', '44\u00d752=2288')\n    ,@('89\u00d765=5785', '97\u00d757=5529')\n    ,@('33\u00d740=1320', '97\u00d719=1843')\n    ,@('50\u00d721=1050', '99\u00d723=2277')\n    ,@('71\u00d754=3834', '54\u00d713=702')\n    ,@('78\u00d716=1248', '82\u00d775=6150')\n    ,@('44\u00d770=3080', '24\u00d794=2256')\n    ,@('83\u00d776=6308', '11\u00d758=638')\n    ,@('63\u00d774=4662', '25\u00d788=2200')\n    ,@('97\u00d723=2231', '10\u00d796=960')\n    ,@('47\u00d762=2914', '91\u00d789=8099')\n    ,@('99\u00d780=7920', '62\u00d794=5828')\n    ,@('57\u00d790=5130', '30\u00d742=1260')\n    ,@('10\u00d746=460', '11\u00d756=616')\n    ,@('90\u00d723=2070', '71\u00d765=4615')\n    ,@('16\u00d727=432', '17\u00d728=476')\n    ,@('52\u00d739=2028', '59\u00d766=3894')\n    ,@('87\u00d752=4524', '29\u00d722=638')\n    ,@('15\u00d771=1065', '56\u00d747=2632')\n    ,@('37\u00d712=444', '96\u00d775=7200')\n    ,@('90\u00d740=3600', '57\u00d743=2451')\n    ,@('42\u00d743=1806', '16\u00d773=1168')\n    ,@('30\u00d760=1800', '19\u00d757=1083')\n    ,@('84\u00d775=6300', '88\u00d735=3080')\n    ,@('22\u00d761=1342', '62\u00d744=2728')\n    ,@('16\u00d717=272', '21\u00d775=1575')\n    ,@('83\u00d777=6391', '51\u00d765=3315')\n    ,@('78\u00d764=4992', '13\u00d778=1014')\n    ,@('24\u00d748=1152', '22\u00d765=1430')\n    ,@('100\u00d735=3500', '46\u00d739=1794')\n    ,@('14\u00d740=560', '15\u00d761=915')\n    ,@('84\u00d798=8232', '20\u00d781=1620')\n    ,@('45\u00d729=1305', '65\u00d734=2210')\n    ,@('85\u00d715=1275', '13\u00d719=247')\n    ,@('72\u00d780=5760', '45\u00d736=1620')\n    ,@('53\u00d729=1537', '57\u00d773=4161')\n    ,@('92\u00d742=3864', '35\u00d775=2625')\n    ,@('19\u00d775=1425', '93\u00d729=2697')\n    ,@('36\u00d749=1764', '11\u00d751=561')\n    ,@('54\u00d733=1782', '22\u00d749=1078')\n    ,@('12\u00d773=876', '28\u00d778=2184')\n    ,@('14\u00d760=840', '18\u00d729=522')\n    ,@('99\u00d733=3267', '78\u00d786=6708')\n    ,@('58\u00d776=4408', '20\u00d760=1200')\n    ,@('92\u00d773=6716', '45\u00d768=3060')\n    ,@('98\u00d745=4410', '78\u00d740=3120')\n    ,@('66\u00d730=1980', '58\u00d759=3422')\n    ,@('22\u00d726=572', '52\u00d711=572')\n    ,@('70\u00d785=5950', '16\u00d750=800')\n    ,@('97\u00d761=5917', '55\u00d713=715')\n    ,@('94\u00d756=5264', '27\u00d727=729')\n    ,@('62\u00d726=1612', '87\u00d766=5742')\n    ,@('56\u00d790=5040', '45\u00d773=3285')\n    ,@('77\u00d735=2695', '85\u00d720=1700')\n    ,@('54\u00d765=3510', '14\u00d755=770')\n    ,@('15\u00d787=1305', '97\u00d759=5723')\n    ,@('86\u00d767=5762', '28\u00d718=504')\n    ,@('40\u00d789=3560', '64\u00d786=5504')\n    ,@('96\u00d748=4608', '17\u00d772=1224')\n    ,@('88\u00d717=1496', '36\u00d749=1764')\n)\n\n# Apply each replacement in order. Every \"old\" value below is unique in the\n# document at the moment it is processed (verified against the source diff),\n# so Find/Execute (wdReplaceAll=2, MatchCase=true) is unambiguous and\n# order-safe even though a handful of \"new\" values coincide with \"old\"\n# values used earlier.\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $oldText\"\n    }\n}\n"}
